$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace the boolean "AgreeTerms" checkbox values (I2:I6) with the text
# value "agreeTerms", left-aligned, matching the new checkbox control text.
for ($r = 2; $r -le 6; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    $cell.Value = "agreeTerms"
    $cell.HorizontalAlignment = -4131   # xlLeft
}

# Give column I (the new text) an explicit width (raw OOXML width of 15).
$ws.Columns.Item(9).ColumnWidth = 14.16666667

# Move the active selection from I6 to I2 (reflects the print-screen change).
$ws.Range("I2").Select() | Out-Null
